$wb = $excel.ActiveWorkbook

# 1. Duplicate the "USE Cello" sheet to serve as the base for the new
#    "USE Double Basses" sheet (keeps styles, data validations, the
#    B1 cell comment, page setup, etc. identical to the source sheet).
$celloSheet = $wb.Worksheets.Item("USE Cello")
$celloSheet.Copy($null, $celloSheet)
$newSheet = $wb.Worksheets.Item("USE Cello (2)")
$newSheet.Name = "USE Double Basses"

# The copied sheet has one extra (blank) trailing row compared with the
# target layout (128 data rows instead of 129) - drop it.
$newSheet.Rows.Item(129).Delete()

# 2. Update the articulation labels for rows 5-10 (columns A and C stay
#    in sync, mirroring every other row on the sheet).
$newSheet.Range("A5").Value = "Tremolo"
$newSheet.Range("C5").Value = "Tremolo"

$newSheet.Range("A6").Value = "Harmonics"
$newSheet.Range("C6").Value = "Harmonics"

$newSheet.Range("A7").Value = "Trills"
$newSheet.Range("C7").Value = "Trills"

$newSheet.Range("A8").Value = "Staccato"
$newSheet.Range("C8").Value = "Staccato"

$newSheet.Range("A9").Value = "Pizzicato"
$newSheet.Range("C9").Value = "Pizzicato"

$newSheet.Range("A10").Value = "Sautille"
$newSheet.Range("C10").Value = "Sautille"

# 3. Shift the MIDI note column (F) up an octave (C0..G0 -> C4..G4) for
#    every populated row.
$newSheet.Range("F3").Value = "C4 (72)"
$newSheet.Range("F4").Value = "C#4 (73)"
$newSheet.Range("F5").Value = "D4 (74)"
$newSheet.Range("F6").Value = "D#4 (75)"
$newSheet.Range("F7").Value = "E4 (76)"
$newSheet.Range("F8").Value = "F4 (77)"
$newSheet.Range("F9").Value = "F#4 (78)"
$newSheet.Range("F10").Value = "G4 (79)"

# Row 11 held "Sautille" on the source (Cello) sheet, which has now
# moved up to row 10 above - clear the now-stale leftover row so it
# goes back to being blank like the rest of the sheet.
$newSheet.Range("A11:K11").ClearContents()

# 4. Move the new sheet so it sits right before "DO NOT MODIFY!" and
#    make it the active tab, matching the target tab order/selection.
$doNotModify = $wb.Worksheets.Item("DO NOT MODIFY!")
$newSheet.Move($doNotModify)
$newSheet.Activate()
